# 自动更新Excel文件 - 2026-01-04 23:13:20
# Decrement the "剩余" (remaining days) column E by 1 for each data row.
# When the remaining days counter reaches 1 (i.e. would become 0), the
# cycle restarts: E is reset to 7 and the start date in column F is
# advanced by 7 days (new restock cycle).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $eVal = $eCell.Value()

    if ($eVal -eq $null) {
        continue
    }

    $fCell = $ws.Cells.Item($r, 6)
    $fVal = $fCell.Value()

    # Rows whose start date is not a well-formed 8-digit yyyymmdd value
    # are left completely untouched (matches source data quirks, e.g. a
    # mistyped 9-digit date that the original updater could not parse).
    if ($fVal -eq $null) {
        continue
    }
    $fStr = [string][int]$fVal
    if ($fStr.Length -ne 8) {
        continue
    }

    if ($eVal -eq 1) {
        # Cycle complete: reset remaining-days counter to 7 and roll the
        # start date forward by 7 days (new restock cycle).
        $y = [int]$fStr.Substring(0, 4)
        $m = [int]$fStr.Substring(4, 2)
        $d = [int]$fStr.Substring(6, 2)

        $dt = Get-Date -Year $y -Month $m -Day $d -Hour 0 -Minute 0 -Second 0
        $dt = $dt.AddDays(7)

        $newF = [int]$dt.ToString("yyyyMMdd")

        $eCell.Value = 7
        $fCell.Value = $newF
    } else {
        $eCell.Value = $eVal - 1
    }
}
